$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra data rows (5-10), leaving header + 3 data rows.
$ws.Rows("5:10").Delete()

# Update row 2
$ws.Range("A2").Value = "file1424-1542116667062.wav"
$ws.Range("B2").Value = 6741
$ws.Range("C2").Value = 48000
$ws.Range("F2").Value = "Kiến Đen uống rượu la đà; Bao nhiêu kiến Gió bay ra chia phần…"

# Update row 3
$ws.Range("A3").Value = "file2608-1542116667082.wav"
$ws.Range("B3").Value = 6570
$ws.Range("C3").Value = 48000
$ws.Range("F3").Value = "Cầm hương kiến Đất bạc đầu; Khóc than kiến Cánh khoác màu áo tang"

# Update row 4
$ws.Range("A4").Value = "file3615-1542116667100.wav"
$ws.Range("B4").Value = 8277
$ws.Range("C4").Value = 48000
$ws.Range("F4").Value = "Sân khấu ở trên không; Giữa vòm trời lá biếc; Trên cành những nhạc công; Cùng thổi kèn náo nhiệt "

# Adjust column widths: B:D -> 10, E -> 30, F -> 50 (was B:F all 50)
# Excel's ColumnWidth setter rounds to pixel granularity before storing the
# OOXML "width" (character-unit) value, so request values just under the
# next pixel boundary to land exactly on 10 / 30 / 50 in the saved file.
$ws.Columns("B:D").ColumnWidth = 9.166666666666666
$ws.Columns("E:E").ColumnWidth = 29.166666666666668
$ws.Columns("F:F").ColumnWidth = 49.166666666666664

Write-Output "done"
